$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on Price (D) cells before assigning their new
# values, so Excel keeps them as text instead of auto-converting numeric-
# looking strings (e.g. "0.9993", "1.000") into floating point numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.526.70"
$ws.Range("E2").Value = "  +1.98%  "

$ws.Range("D3").Value = "1.665.69"
$ws.Range("E3").Value = "  +0.93%  "

$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "237.56"
$ws.Range("E5").Value = "  +0.24%  "

$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").Value = "0.4798"
$ws.Range("E7").Value = "  -0.66%  "

$ws.Range("D8").Value = "0.2626"
$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("D9").Value = "0.06165"
$ws.Range("E9").Value = "  +2.29%  "

$ws.Range("D10").Value = "0.07084"
$ws.Range("E10").Value = "  -1.60%  "

$ws.Range("D11").Value = "1.664.67"
$ws.Range("E11").Value = "  +0.96%  "

$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("D13").Value = "0.5870"
$ws.Range("E13").Value = "  -5.47%  "

$ws.Range("D14").Value = "4.363"
$ws.Range("E14").Value = "  -4.75%  "

$ws.Range("D15").Value = "74.94"
$ws.Range("E15").Value = "  +2.67%  "

$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.04%  "

$ws.Range("D17").Value = "0.9998"
$ws.Range("E17").Value = "  +0.10%  "

$ws.Range("D18").Value = "25.529.75"
$ws.Range("E18").Value = "  +2.05%  "

$ws.Range("D19").Value = "0.000006751"
$ws.Range("E19").Value = "  +1.85%  "

$ws.Range("D20").Value = "11.45"
$ws.Range("E20").Value = "  -0.54%  "

$ws.Range("D21").Value = "1.877.43"
$ws.Range("E21").Value = "  +1.07%  "

$ws.Range("D22").Value = "4.409"
$ws.Range("E22").Value = "  -3.14%  "

$ws.Range("D23").Value = "8.724"
$ws.Range("E23").Value = "  +1.21%  "

$ws.Range("D24").Value = "5.275"
$ws.Range("E24").Value = "  -0.59%  "

$ws.Range("D25").Value = "135.64"
$ws.Range("E25").Value = "  +2.70%  "

$ws.Range("D26").Value = "15.01"
$ws.Range("E26").Value = "  +0.47%  "

$ws.Range("E27").Value = "  -0.40%  "

$ws.Range("D28").Value = "104.80"
$ws.Range("E28").Value = "  +1.64%  "

$ws.Range("D29").Value = "1.709"
$ws.Range("E29").Value = "  +1.84%  "

$ws.Range("D30").Value = "3.968"
$ws.Range("E30").Value = "  +5.49%  "

$ws.Range("D31").Value = "0.07760"
$ws.Range("E31").Value = "  -1.66%  "

$ws.Range("D32").Value = "3.636"
$ws.Range("E32").Value = "  +1.08%  "

$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("D34").Value = "0.04209"
$ws.Range("E34").Value = "  -7.94%  "

$ws.Range("D35").Value = "2.599"
$ws.Range("E35").Value = "  +0.26%  "

$ws.Range("D36").Value = "0.6095"
$ws.Range("E36").Value = "  +5.25%  "

$ws.Range("D37").Value = "0.9467"
$ws.Range("E37").Value = "  +1.05%  "

$ws.Range("E38").Value = "  -0.12%  "

$ws.Range("D39").Value = "0.8620"
$ws.Range("E39").Value = "  +2.01%  "

$ws.Range("D40").Value = "1.001"
$ws.Range("E40").Value = "  +0.13%  "

$ws.Range("D41").Value = "1.848"
$ws.Range("E41").Value = "  +0.99%  "

$ws.Range("D42").Value = "0.01461"
$ws.Range("E42").Value = "  -6.44%  "

$ws.Range("D43").Value = "97.15"
$ws.Range("E43").Value = "  -1.22%  "

$ws.Range("D44").Value = "0.3754"
$ws.Range("E44").Value = "  +0.65%  "

$ws.Range("D45").Value = "4.836"
$ws.Range("E45").Value = "  +1.11%  "

$ws.Range("E46").Value = "  -2.02%  "

$ws.Range("D47").Value = "6.192"
$ws.Range("E47").Value = "  +0.68%  "

$ws.Range("D48").Value = "0.05264"
$ws.Range("E48").Value = "  +1.42%  "

$ws.Range("D49").Value = "29.71"
$ws.Range("E49").Value = "  -0.39%  "

$ws.Range("B50").Value = "TrueUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
$ws.Range("D50").Value = "1.002"
$ws.Range("E50").Value = "  +0.11%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.355"
$ws.Range("E51").Value = "  +2.44%  "
